# Update gh-pages to output generated at 456a3b4
# Applies updated "want to go" counts, a renamed event title, and a
# refreshed cover image URL to both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$newTitle = "【会员购严选】苏州·back to ACG端阳嘉年华动漫国潮文化节"
$newCover = "//i0.hdslb.com/bfs/openplatform/202405/vPI9YxcW1715674161718.jpeg"

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value = 130
$ws1.Range("F4").Value = 2081
$ws1.Range("F5").Value = 362
$ws1.Range("F6").Value = 629

$ws1.Range("C9").Value = $newTitle
$ws1.Range("F9").Value = 10690
$ws1.Range("I9").Value = $newCover

$ws1.Range("F14").Value = 417
$ws1.Range("F15").Value = 7542
$ws1.Range("F17").Value = 721
$ws1.Range("F18").Value = 259

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value = 130
$ws4.Range("F4").Value = 2081
$ws4.Range("F5").Value = 362
$ws4.Range("F6").Value = 629

$ws4.Range("C12").Value = $newTitle
$ws4.Range("F12").Value = 10690
$ws4.Range("I12").Value = $newCover

$ws4.Range("F17").Value = 417
$ws4.Range("F18").Value = 7542
$ws4.Range("F20").Value = 721
$ws4.Range("F21").Value = 259
